# "2022-03-11 Table in PlutoWeb.xlsx" - beta 24.12.2022
# Adds two new "areaAccount" action elements (ac_button_verify_account,
# ac_verification_hint) and inserts blank separator rows above each
# "Area | Action Elements" sub-block, shifting everything below down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank separator row above the "areaCreateAccount" action-elements block
# (was row 15, now becomes row 16 and onward).
$ws.Rows("15:15").Insert() | Out-Null
$ws.Rows("15:15").RowHeight = 15

# Blank separator row above the "areaAccount" action-elements block
# (was row 20, now row 21 after the previous insert).
$ws.Rows("21:21").Insert() | Out-Null
$ws.Rows("21:21").RowHeight = 15

# Two new rows for "areaAccount" action elements plus a trailing blank
# separator row above the "areaMessages" block (inserted right after the
# former "ac_email" row, which is now row 24).
$ws.Rows("25:27").Insert() | Out-Null
$ws.Rows("25:27").RowHeight = 15

$ws.Range("C25").Value = "ac_button_verify_account"
$ws.Range("C26").Value = "ac_verification_hint"

# Reflect the author's new viewport/selection position.
$ws.Range("C19").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
